## Configured Mail with attachments
## Appends four new log rows (S.No 7-10) to the "sriram" timesheet sheet,
## introduces the built-in "Bad" cell style for the newest "started" status,
## and moves the active selection to the newly entered cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sriram")
$ws.Activate()

# --- Row 13 : reuse the formatting already used by row 9 (A wrap, B date,
#              C/D wrap, E "Good") then overwrite the values. ---
$ws.Range("A9:E9").Copy()
$ws.Range("A13:E13").PasteSpecial(-4122)
$ws.Range("A13").Value = 7
$ws.Range("B13").Value = 43124
$ws.Range("C13").Value = "11:30 to 17:00"
$ws.Range("D13").Value = "Configured Mail using outlook mail"
$ws.Range("E13").Value = "completed"

# --- Row 14 : same as row 13 except S.No (A14) keeps no explicit style.
#              Project task (D) was typed before the time (C). ---
$ws.Range("A9:E9").Copy()
$ws.Range("A14:E14").PasteSpecial(-4122)
$ws.Range("A14").Value = 8
$ws.Range("A14").ClearFormats()
$ws.Range("B14").Value = 43125
$ws.Range("D14").Value = "Template for mail"
$ws.Range("C14").Value = "13:00 to 17:00"
$ws.Range("E14").Value = "completed"

# --- Row 15 : identical pattern to row 13; again D typed before C. ---
$ws.Range("A9:E9").Copy()
$ws.Range("A15:E15").PasteSpecial(-4122)
$ws.Range("A15").Value = 9
$ws.Range("B15").Value = 43129
$ws.Range("D15").Value = "Configured mail with attachments"
$ws.Range("C15").Value = "13:00 to 14:00"
$ws.Range("E15").Value = "completed"

# --- Row 16 : final row, no date cell, status flagged with the "Bad" style. ---
$ws.Range("C10:D10").Copy()
$ws.Range("C16:D16").PasteSpecial(-4122)
$ws.Range("A16").Value = 10
$ws.Range("C16").Value = "14:00 to 17:00"
$ws.Range("D16").Value = "Filtering Data from excel"
$ws.Range("E16").Value = "started"
$ws.Range("E16").Style = "Bad"

$excel.CutCopyMode = $false

# Move the visible selection onto the freshly typed last cell, matching
# where the author's cursor ended up after the edit.
$ws.Range("E16").Select()
